# feat: add 2022-Q4 data
#
# Before: sheet1 "总计" (summary) + sheet2 "2020-Q4" (fund holdings detail)
# After:  sheet1 "总计" (summary, now with a 2022-Q4 row too) +
#         sheet2 "2022-Q4" (new fund holdings detail) +
#         sheet3 "2020-Q4" (old fund holdings detail, unchanged, relocated)

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# Helper: write $value into $range as TEXT (not auto-coerced to a
# number) while leaving the cell's style untouched (no "@" / quote
# prefix residue) by pasting the number-format from a pristine,
# default-styled cell on top afterwards.
# ---------------------------------------------------------------------
function Set-TextValue($range, $value, $blank) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $blank.Copy()
    $range.PasteSpecial(-4122)
}

# Writes one data row of the fund-holdings sheets (columns A..H).
function Set-FundRow($ws, $r, $idx, $code, $name, $scale, $pos, $ratio, $mv, $rank, $blank) {
    $ws.Cells.Item($r, 1).Value = $idx
    Set-TextValue $ws.Cells.Item($r, 2) $code $blank
    Set-TextValue $ws.Cells.Item($r, 3) $name $blank
    Set-TextValue $ws.Cells.Item($r, 4) $scale $blank
    Set-TextValue $ws.Cells.Item($r, 5) $pos $blank
    Set-TextValue $ws.Cells.Item($r, 6) $ratio $blank
    Set-TextValue $ws.Cells.Item($r, 7) $mv $blank
    $ws.Cells.Item($r, 8).Value = $rank
}

# =======================================================================
# 1. Re-organise the sheets: the existing "2020-Q4" sheet becomes
#    "2022-Q4" (it keeps sheetId 2 / rId2), and a brand-new sheet named
#    "2020-Q4" is added right after it (sheetId 3 / rId3) to hold the
#    original data that used to live there.
# =======================================================================
$ws2.Name = "2022-Q4"
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "2020-Q4"

# =======================================================================
# 2. Populate the new "2020-Q4" sheet (sheet3) with the data that used
#    to be on the "2020-Q4" sheet before the edit. This sheet uses the
#    OTHER header style (bold/bordered, but without the explicit
#    theme-color/size that sheet1's style carries) - reproduce it via
#    the Font/Border object model, which resolves to that exact,
#    already-present style rather than minting a new one.
# =======================================================================
$blank3 = $ws3.Range("Z100")

$ws3.Range("B1").Value = "基金代码"
$ws3.Range("C1").Value = "基金名称"
$ws3.Range("D1").Value = "基金金额"
$ws3.Range("E1").Value = "股票总仓位"
$ws3.Range("F1").Value = "仓位占比"
$ws3.Range("G1").Value = "持有市值(亿元)"
$ws3.Range("H1").Value = "仓位排名"

$hdr3 = $ws3.Range("B1:H1")
$hdr3.Font.Bold = $true
$hdr3.HorizontalAlignment = -4108
$hdr3.VerticalAlignment = -4160
$hdr3.Borders.LineStyle = 1

Set-FundRow $ws3 2 0 "005702" "恒生前海港股通高股息低波动指数" "0.31" "94.32" "2.64" "0.0082" 4 $blank3

$a2_3 = $ws3.Range("A2")
$a2_3.Font.Bold = $true
$a2_3.HorizontalAlignment = -4108
$a2_3.VerticalAlignment = -4160
$a2_3.Borders.LineStyle = 1

# =======================================================================
# 3. Replace the "2022-Q4" sheet (sheet2, still holding the old data)
#    with the brand-new fund-holdings data for 2022-Q4.
# =======================================================================
$ws2.Cells.Clear()
$blank2 = $ws2.Range("Z100")

$ws2.Range("B1").Value = "基金代码"
$ws2.Range("C1").Value = "基金名称"
$ws2.Range("D1").Value = "基金规模"
$ws2.Range("E1").Value = "股票总仓位"
$ws2.Range("F1").Value = "仓位占比"
$ws2.Range("G1").Value = "持有市值(亿元)"
$ws2.Range("H1").Value = "仓位排名"
$ws1.Range("B1:D1").Copy()
$ws2.Range("B1:H1").PasteSpecial(-4122)

Set-FundRow $ws2 2 0 "007368" "浙商沪港深精选混合A"             "6.05" "92.12" "5.12" "0.3098" 7  $blank2
Set-FundRow $ws2 3 1 "010381" "浙商智选价值混合A"               "7.03" "91.16" "3.36" "0.2362" 10 $blank2
Set-FundRow $ws2 4 2 "010382" "浙商智选价值混合C"               "6.53" "91.16" "3.36" "0.2194" 10 $blank2
Set-FundRow $ws2 5 3 "007369" "浙商沪港深精选混合C"             "1.38" "92.12" "5.12" "0.0707" 7  $blank2
Set-FundRow $ws2 6 4 "016518" "华泰紫金创新成长混合C"           "0.81" "94.10" "4.61" "0.0373" 10 $blank2
Set-FundRow $ws2 7 5 "016517" "华泰紫金创新成长混合A"           "0.49" "94.10" "4.61" "0.0226" 10 $blank2
Set-FundRow $ws2 8 6 "009569" "浙商智多宝稳健一年持有期混合C"   "1.02" "26.91" "1.18" "0.0120" 6  $blank2
Set-FundRow $ws2 9 7 "009568" "浙商智多宝稳健一年持有期混合A"   "0.98" "26.91" "1.18" "0.0116" 6  $blank2

# Column A (the row-index column) carries the same header style as
# sheet1's index column - copy it across the whole A2:A9 block at once.
$ws1.Range("A2").Copy()
$ws2.Range("A2:A9").PasteSpecial(-4122)
$ws2.Range("A2").Value = 0
$ws2.Range("A3").Value = 1
$ws2.Range("A4").Value = 2
$ws2.Range("A5").Value = 3
$ws2.Range("A6").Value = 4
$ws2.Range("A7").Value = 5
$ws2.Range("A8").Value = 6
$ws2.Range("A9").Value = 7

# Match sheet1's page margins (the new sheet was authored like sheet1,
# not with Excel's 0.7/0.75/0.3 worksheet defaults).
$ws2.PageSetup.LeftMargin = 54
$ws2.PageSetup.RightMargin = 54
$ws2.PageSetup.TopMargin = 72
$ws2.PageSetup.BottomMargin = 72
$ws2.PageSetup.HeaderMargin = 36
$ws2.PageSetup.FooterMargin = 36

# =======================================================================
# 4. Update the "总计" summary sheet: insert the 2022-Q4 total as the
#    new row 2 and push the 2020-Q4 total down to row 3.
# =======================================================================
$ws1.Range("A2").Copy($ws1.Range("A3"))
$ws1.Range("A3").Value = 1
$ws1.Range("B3").Value = "2020-Q4"
$ws1.Range("C3").Value = 1
$ws1.Range("D3").Value = 0.01

$ws1.Range("B2").Value = "2022-Q4"
$ws1.Range("C2").Value = 8
$ws1.Range("D2").Value = 0.92

$ws1.Activate()
